# Applies the "made minor changes on self-control" edit to the Self-control
# section of the technical note.
#
# Three small wording/typo fixes in the "Self-control" paragraphs:
#   1. "...processing of a information..." -> "...processing of an information..."
#   2. "...briefly presented o the participant..." -> "...briefly presented to the participant..."
#   3. "...(beneficial condition of BC), and the inverse situation (control condition or CC)..."
#      -> "...(beneficial condition or BC), and the inverse situation (detrimental condition or DC)..."

$d = $word.ActiveDocument

# 1) "of a information" -> "of an information"
$d.Content.Find.Execute(
    "prioritize the attentional processing of a information",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "prioritize the attentional processing of an information",
    2)

# 2) "presented o the participant" -> "presented to the participant"
$d.Content.Find.Execute(
    "briefly presented o the participant",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "briefly presented to the participant",
    2)

# 3) "condition of BC" -> "condition or BC" and "control condition or CC" -> "detrimental condition or DC"
$d.Content.Find.Execute(
    "beneficial condition of BC), and the inverse situation (control condition or CC)",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "beneficial condition or BC), and the inverse situation (detrimental condition or DC)",
    2)
